# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the row for
# c476e3d3-e3a6-42fe-a53c-fc9112a467ad.md has now been handed off for
# translation, so its status moves from "In Translation" to "Ready for
# handoff", its priority switches to machine translation ("mt"), and the
# handoff timestamps are refreshed to the moment the handoff xliff files
# were generated. (Status text got longer, so the Status/locale columns
# also widen a touch to fit it.)

$wb = $excel.ActiveWorkbook

# --- Overview sheet ----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-24 18:21:10"
$wsOverview.Columns.Item(5).ColumnWidth = 16.85
$wsOverview.Columns.Item(6).ColumnWidth = 16.85

# --- zh-cn sheet ---------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-24 18:20:58"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.85

# --- de-de sheet ---------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-24 18:21:10"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.85
